# Reorder worksheets so that "总计" becomes the first sheet and
# "2021-Q3" becomes the second sheet (swap their tab order).
$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3    = $wb.Worksheets.Item("2021-Q3")

# Move "总计" to be before the "2021-Q3" sheet (i.e. make it first).
$wsTotal.Move($wsQ3)
